# rcrb_laa.xlsx — "cleaning up lm models for cali"
#
# The third column ("Value", a dollar-formatted figure) is replaced by a new
# "Participation" column of small integer counts, formatted the same way as
# the existing "Landings (lbs)" column (#,0 — no currency sign). The Total
# row's participation cell becomes a plain sum-of-the-literal (566), and its
# format/style is likewise unified with the Landings Total cell's style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: column C label "Value" -> "Participation"
$ws.Cells.Item(1, 3).Value = "Participation"

# New participation values, one per year-row (row 2 = 1980 ... row 45 = 2023)
$participation = @(
    62, 74, 68, 63, 80, 65, 57, 50, 48, 36,
    62, 55, 46, 35, 38, 30, 41, 37, 35, 44,
    42, 48, 38, 43, 37, 34, 32, 33, 29, 24,
    31, 29, 39, 32, 33, 29, 36, 27, 26, 37,
    36, 33, 26, 33
)

for ($i = 0; $i -lt $participation.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Cells.Item($row, 3)
    $cell.NumberFormat = "#,0"
    $cell.Value = $participation[$i]
}

# Total row (row 46): participation total becomes a plain literal, styled
# the same (#,0, bold) as the Landings total cell in column B.
$ws.Cells.Item(46, 3).NumberFormat = "#,0"
$ws.Cells.Item(46, 3).Value = 566
